# Rename sample sheet header columns (SAMPLESHEET_TO_BARCODE)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Sample Name"
$ws.Range("B1").Value = "5' Barcode Sequence"
$ws.Range("C1").Value = "3' Barcode Sequence"

# Mirror the residual selection state saved in the workbook after editing
$ws.Range("A1:G3").Select()
